# Refresh cached market-price / profit figures on the per-class Leve sheets.
# Mirrors a scheduled scraper run: only the price/profit columns (H-N) on a
# handful of existing rows change; no rows/columns are inserted or removed.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19: H19, I19, J19, K19, L19, M19, N19
$ws.Range("H19").Value = 2574.2
$ws.Range("I19").Value = 4560
$ws.Range("J19").Value = 588.4
$ws.Range("K19").Value = 4560
$ws.Range("L19").Value = 588.4
$ws.Range("M19").Value = -4385
$ws.Range("N19").Value = -938.4
# Row 40: H40, I40, J40, K40, L40, M40, N40
$ws.Range("H40").Value = 1977949.1
$ws.Range("I40").Value = 1035.4445
$ws.Range("J40").Value = 3346581.8
$ws.Range("K40").Value = 1035.4445
$ws.Range("L40").Value = 3346581.8
$ws.Range("M40").Value = -860.4445000000001
$ws.Range("N40").Value = -3346931.8
# Row 111: H111, J111, L111, N111
$ws.Range("H111").Value = 3642
$ws.Range("J111").Value = 1643.6666
$ws.Range("L111").Value = 4930.9998
$ws.Range("N111").Value = -11064.9998
# Row 127: H127, J127, L127, N127
$ws.Range("H127").Value = 1142.8572
$ws.Range("J127").Value = 1315.375
$ws.Range("L127").Value = 3946.125
$ws.Range("N127").Value = -13866.125
# Row 129: H129, J129, L129, N129
$ws.Range("H129").Value = 233624.81
$ws.Range("J129").Value = 278986.38
$ws.Range("L129").Value = 836959.14
$ws.Range("N129").Value = -846959.14
# Row 137: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 22896.66
$ws.Range("I137").Value = 1407.2354
$ws.Range("J137").Value = 79099.766
$ws.Range("K137").Value = 4221.706200000001
$ws.Range("L137").Value = 237299.298
$ws.Range("M137").Value = -1671.706200000001
$ws.Range("N137").Value = -242399.298
# Row 138: H138, J138, L138, N138
$ws.Range("H138").Value = 2373.7866
$ws.Range("J138").Value = 2476.9
$ws.Range("L138").Value = 7430.700000000001
$ws.Range("N138").Value = -17710.7

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32, I32, J32, K32, L32, M32, N32
$ws.Range("H32").Value = 23463.092
$ws.Range("I32").Value = 26087.396
$ws.Range("J32").Value = 2468.6667
$ws.Range("K32").Value = 26087.396
$ws.Range("L32").Value = 2468.6667
$ws.Range("M32").Value = -25800.396
$ws.Range("N32").Value = -3042.6667
# Row 45: H45, I45, J45, K45, L45, M45, N45
$ws.Range("H45").Value = 2722.875
$ws.Range("I45").Value = 2456.72
$ws.Range("J45").Value = 3673.4285
$ws.Range("K45").Value = 2456.72
$ws.Range("L45").Value = 3673.4285
$ws.Range("M45").Value = -2079.72
$ws.Range("N45").Value = -4427.4285
# Row 74: H74, I74, K74, M74
$ws.Range("H74").Value = 125001020
$ws.Range("I74").Value = 125001020
$ws.Range("K74").Value = 125001020
$ws.Range("M74").Value = -125000146
# Row 77: H77, I77, K77, M77
$ws.Range("H77").Value = 125001020
$ws.Range("I77").Value = 125001020
$ws.Range("K77").Value = 625005100
$ws.Range("M77").Value = -625000732
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 2961
$ws.Range("I122").Value = 2305.818
$ws.Range("J122").Value = 4402.4
$ws.Range("K122").Value = 6917.454000000001
$ws.Range("L122").Value = 13207.2
$ws.Range("M122").Value = -4467.454000000001
$ws.Range("N122").Value = -18107.2
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 33151.668
$ws.Range("I132").Value = 1792.4048
$ws.Range("K132").Value = 5377.2144
$ws.Range("M132").Value = -2847.2144

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80: H80, I80, J80, K80, L80, M80, N80
$ws.Range("H80").Value = 756.0968
$ws.Range("I80").Value = 649.8461
$ws.Range("J80").Value = 832.8333
$ws.Range("K80").Value = 649.8461
$ws.Range("L80").Value = 832.8333
$ws.Range("M80").Value = 348.1539
$ws.Range("N80").Value = -2828.8333
# Row 83: H83, I83, J83, K83, L83, M83, N83
$ws.Range("H83").Value = 756.0968
$ws.Range("I83").Value = 649.8461
$ws.Range("J83").Value = 832.8333
$ws.Range("K83").Value = 3249.2305
$ws.Range("L83").Value = 4164.1665
$ws.Range("M83").Value = 1742.7695
$ws.Range("N83").Value = -14148.1665
# Row 105: H105, I105, J105, K105, L105, M105, N105
$ws.Range("H105").Value = 2536.7727
$ws.Range("I105").Value = 2693.077
$ws.Range("J105").Value = 2311
$ws.Range("K105").Value = 2693.077
$ws.Range("L105").Value = 2311
$ws.Range("M105").Value = -946.0770000000002
$ws.Range("N105").Value = -5805
# Row 107: H107, I107, K107, M107
$ws.Range("H107").Value = 933.4400000000001
$ws.Range("I107").Value = 750.9375
$ws.Range("K107").Value = 750.9375
$ws.Range("M107").Value = 1169.0625
# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 37904.566
$ws.Range("I134").Value = 43216.81
$ws.Range("J134").Value = 3375
$ws.Range("K134").Value = 129650.43
$ws.Range("L134").Value = 10125
$ws.Range("M134").Value = -127115.43
$ws.Range("N134").Value = -15195

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 20536.105
$ws.Range("I31").Value = 27961.334
$ws.Range("J31").Value = 7807.143
$ws.Range("K31").Value = 27961.334
$ws.Range("L31").Value = 7807.143
$ws.Range("M31").Value = -27666.334
$ws.Range("N31").Value = -8397.143
# Row 34: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 20536.105
$ws.Range("I34").Value = 27961.334
$ws.Range("J34").Value = 7807.143
$ws.Range("K34").Value = 27961.334
$ws.Range("L34").Value = 7807.143
$ws.Range("M34").Value = -27759.334
$ws.Range("N34").Value = -8211.143
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 18862.773
$ws.Range("I132").Value = 21828.24
$ws.Range("K132").Value = 65484.72
$ws.Range("M132").Value = -62954.72
# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 920.2143
$ws.Range("I134").Value = 762.0909
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 2286.2727
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = 248.7273
$ws.Range("N134").Value = -9570

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3: H3, J3, L3, N3
$ws.Range("H3").Value = 2895.8333
$ws.Range("J3").Value = 4646
$ws.Range("L3").Value = 13938
$ws.Range("N3").Value = -14162
# Row 5: H5, I5, J5, K5, L5, M5, N5
$ws.Range("H5").Value = 1377
$ws.Range("I5").Value = 1121.9
$ws.Range("J5").Value = 2652.5
$ws.Range("K5").Value = 3365.7
$ws.Range("L5").Value = 7957.5
$ws.Range("M5").Value = -3253.7
$ws.Range("N5").Value = -8181.5
# Row 92: H92, J92, L92, N92
$ws.Range("H92").Value = 12500457
$ws.Range("J92").Value = 514.8333
$ws.Range("L92").Value = 1544.4999
$ws.Range("N92").Value = -4040.4999
# Row 97: H97, I97, K97, M97
$ws.Range("H97").Value = 906.25
$ws.Range("I97").Value = 250
$ws.Range("K97").Value = 750
$ws.Range("M97").Value = -254
# Row 98: H98, J98, L98, N98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 107: H107, J107, L107, N107
$ws.Range("H107").Value = 4533.609
$ws.Range("J107").Value = 239.5
$ws.Range("L107").Value = 718.5
$ws.Range("N107").Value = -4558.5
# Row 114: H114, I114, J114, K114, L114, M114, N114
$ws.Range("H114").Value = 3094.7778
$ws.Range("I114").Value = 1990.8
$ws.Range("J114").Value = 4474.75
$ws.Range("K114").Value = 5972.4
$ws.Range("L114").Value = 13424.25
$ws.Range("M114").Value = -2718.4
$ws.Range("N114").Value = -19932.25
# Row 129: H129, I129, J129, K129, L129, M129, N129
$ws.Range("H129").Value = 263865.94
$ws.Range("I129").Value = 450.55554
$ws.Range("J129").Value = 500939.8
$ws.Range("K129").Value = 1351.66662
$ws.Range("L129").Value = 1502819.4
$ws.Range("M129").Value = 3648.33338
$ws.Range("N129").Value = -1512819.4
# Row 131: H131, I131, J131, K131, L131, M131, N131
$ws.Range("H131").Value = 751.1818
$ws.Range("I131").Value = 313
$ws.Range("J131").Value = 800.4157
$ws.Range("K131").Value = 939
$ws.Range("L131").Value = 2401.2471
$ws.Range("M131").Value = 4101
$ws.Range("N131").Value = -12481.2471
# Row 133: H133, I133, J133, K133, L133, M133, N133
$ws.Range("H133").Value = 4090
$ws.Range("I133").Value = 2180
$ws.Range("J133").Value = 6000
$ws.Range("K133").Value = 6540
$ws.Range("L133").Value = 18000
$ws.Range("M133").Value = -1480
$ws.Range("N133").Value = -28120
# Row 135: H135, I135, J135, K135, L135, M135, N135
$ws.Range("H135").Value = 1377
$ws.Range("I135").Value = 1121.9
$ws.Range("J135").Value = 2652.5
$ws.Range("K135").Value = 10097.1
$ws.Range("L135").Value = 23872.5
$ws.Range("M135").Value = -7562.1
$ws.Range("N135").Value = -28942.5

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16: H16, I16, K16, M16
$ws.Range("H16").Value = 499.66666
$ws.Range("I16").Value = 499.66666
$ws.Range("K16").Value = 499.66666
$ws.Range("M16").Value = -329.66666
# Row 22: H22, I22, J22, K22, L22, M22, N22
$ws.Range("H22").Value = 2447.5
$ws.Range("I22").Value = 2845.25
$ws.Range("J22").Value = 2182.3333
$ws.Range("K22").Value = 2845.25
$ws.Range("L22").Value = 2182.3333
$ws.Range("M22").Value = -2550.25
$ws.Range("N22").Value = -2772.3333
# Row 27: H27, I27, J27, K27, L27, M27, N27
$ws.Range("H27").Value = 2447.5
$ws.Range("I27").Value = 2845.25
$ws.Range("J27").Value = 2182.3333
$ws.Range("K27").Value = 2845.25
$ws.Range("L27").Value = 2182.3333
$ws.Range("M27").Value = -2738.25
$ws.Range("N27").Value = -2396.3333
# Row 68: H68, I68, J68, K68, L68, M68, N68
$ws.Range("H68").Value = 2417.5881
$ws.Range("I68").Value = 2391.8333
$ws.Range("J68").Value = 2479.4
$ws.Range("K68").Value = 2391.8333
$ws.Range("L68").Value = 2479.4
$ws.Range("M68").Value = -1642.8333
$ws.Range("N68").Value = -3977.4
# Row 71: H71, I71, J71, K71, L71, M71, N71
$ws.Range("H71").Value = 2417.5881
$ws.Range("I71").Value = 2391.8333
$ws.Range("J71").Value = 2479.4
$ws.Range("K71").Value = 11959.1665
$ws.Range("L71").Value = 12397
$ws.Range("M71").Value = -8215.166499999999
$ws.Range("N71").Value = -19885
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1296.6216
$ws.Range("I132").Value = 969
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 2907
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -377
$ws.Range("N132").Value = -17058.5
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 32752.562
$ws.Range("I136").Value = 36895.43
$ws.Range("K136").Value = 110686.29
$ws.Range("M136").Value = -108136.29

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 27028418
$ws.Range("I136").Value = 29412970
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 88238910
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -88236360
$ws.Range("N136").Value = -15600
